$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.019999999999999
$ws.Cells.Item(2, 3).Value = 1.031485694350901
$ws.Cells.Item(2, 4).Value = 1.034130082505344
$ws.Cells.Item(2, 5).Value = 1.039945419585241
$ws.Cells.Item(2, 6).Value = 1.048339139711061
$ws.Cells.Item(2, 9).Value = 1.027121246648626
$ws.Cells.Item(2, 10).Value = 1.036621158574207
$ws.Cells.Item(2, 11).Value = 1.036930441628195
$ws.Cells.Item(2, 12).Value = 1.042729166280559
$ws.Cells.Item(2, 13).Value = 1.051099260956818
$ws.Cells.Item(2, 14).Value = 1.038093279521491

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.032914841727663
$ws.Cells.Item(3, 4).Value = 1.035472907661368
$ws.Cells.Item(3, 5).Value = 1.041211508479991
$ws.Cells.Item(3, 6).Value = 1.049684684210353
$ws.Cells.Item(3, 9).Value = 1.027107386996032
$ws.Cells.Item(3, 10).Value = 1.037689760018518
$ws.Cells.Item(3, 11).Value = 1.038080968582189
$ws.Cells.Item(3, 12).Value = 1.043804377773365
$ws.Cells.Item(3, 13).Value = 1.052255450554406
$ws.Cells.Item(3, 14).Value = 1.039163398502424

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.033838931196638
$ws.Cells.Item(4, 4).Value = 1.036341473076731
$ws.Cells.Item(4, 5).Value = 1.042030183676856
$ws.Cells.Item(4, 6).Value = 1.050554198088799
$ws.Cells.Item(4, 9).Value = 1.02709594097996
$ws.Cells.Item(4, 10).Value = 1.038380165121777
$ws.Cells.Item(4, 11).Value = 1.038824574530844
$ws.Cells.Item(4, 12).Value = 1.044498993703159
$ws.Cells.Item(4, 13).Value = 1.053001905458458
$ws.Cells.Item(4, 14).Value = 1.03985478406012

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.034227264771045
$ws.Cells.Item(5, 4).Value = 1.036706542765868
$ws.Cells.Item(5, 5).Value = 1.042374222274319
$ws.Cells.Item(5, 6).Value = 1.050919471727514
$ws.Cells.Item(5, 9).Value = 1.027090535141004
$ws.Cells.Item(5, 10).Value = 1.038670163352662
$ws.Cells.Item(5, 11).Value = 1.039136983990695
$ws.Cells.Item(5, 12).Value = 1.044790746023724
$ws.Cells.Item(5, 13).Value = 1.053315316837769
$ws.Cells.Item(5, 14).Value = 1.040145194121756

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.034292458834961
$ws.Cells.Item(6, 4).Value = 1.036767835288093
$ws.Cells.Item(6, 5).Value = 1.042431980210292
$ws.Cells.Item(6, 6).Value = 1.050980787054044
$ws.Cells.Item(6, 9).Value = 1.027089592627047
$ws.Cells.Item(6, 10).Value = 1.038718840890289
$ws.Cells.Item(6, 11).Value = 1.039189427152059
$ws.Cells.Item(6, 12).Value = 1.044839717126425
$ws.Cells.Item(6, 13).Value = 1.053367916726609
$ws.Cells.Item(6, 14).Value = 1.040193940787069

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.033844120725848
$ws.Cells.Item(7, 4).Value = 1.036346351445973
$ws.Cells.Item(7, 5).Value = 1.042034781253274
$ws.Cells.Item(7, 6).Value = 1.050559079949389
$ws.Cells.Item(7, 9).Value = 1.027095871081568
$ws.Cells.Item(7, 10).Value = 1.038384041062604
$ws.Cells.Item(7, 11).Value = 1.038828749752192
$ws.Cells.Item(7, 12).Value = 1.044502893145985
$ws.Cells.Item(7, 13).Value = 1.053006094837627
$ws.Cells.Item(7, 14).Value = 1.039858665505227

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.031968822251985
$ws.Cells.Item(8, 4).Value = 1.034583967961278
$ws.Cells.Item(8, 5).Value = 1.04037342033468
$ws.Cells.Item(8, 6).Value = 1.048794111774693
$ws.Cells.Item(8, 9).Value = 1.027117075349439
$ws.Cells.Item(8, 10).Value = 1.036982517357087
$ws.Cells.Item(8, 11).Value = 1.037319448398184
$ws.Cells.Item(8, 12).Value = 1.04309277253226
$ws.Cells.Item(8, 13).Value = 1.051490349848124
$ws.Cells.Item(8, 14).Value = 1.038455151475292

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.02865897405227
$ws.Cells.Item(9, 4).Value = 1.031475683801637
$ws.Cells.Item(9, 5).Value = 1.037441363213987
$ws.Cells.Item(9, 6).Value = 1.045675086972748
$ws.Cells.Item(9, 9).Value = 1.027135493338906
$ws.Cells.Item(9, 10).Value = 1.03450463607754
$ws.Cells.Item(9, 11).Value = 1.034653091461776
$ws.Cells.Item(9, 12).Value = 1.040599246237562
$ws.Cells.Item(9, 13).Value = 1.04880642305911
$ws.Cells.Item(9, 14).Value = 1.035973751320111

$ws.Cells.Item(10, 2).Value = 1.019999999999999
$ws.Cells.Item(10, 3).Value = 1.026448480574685
$ws.Cells.Item(10, 4).Value = 1.029401378332018
$ws.Cells.Item(10, 5).Value = 1.035483375659143
$ws.Cells.Item(10, 6).Value = 1.043589497278466
$ws.Cells.Item(10, 9).Value = 1.027135073965822
$ws.Cells.Item(10, 10).Value = 1.032846958485051
$ws.Cells.Item(10, 11).Value = 1.032870719864811
$ws.Cells.Item(10, 12).Value = 1.038930821087465
$ws.Cells.Item(10, 13).Value = 1.047008196840262
$ws.Cells.Item(10, 14).Value = 1.034313719635302

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.02549030406946
$ws.Cells.Item(11, 4).Value = 1.028502617583057
$ws.Cells.Item(11, 5).Value = 1.034634712796473
$ws.Cells.Item(11, 6).Value = 1.042684879697908
$ws.Cells.Item(11, 9).Value = 1.027131888730443
$ws.Cells.Item(11, 10).Value = 1.032127748160393
$ws.Cells.Item(11, 11).Value = 1.03209773954214
$ws.Cells.Item(11, 12).Value = 1.03820688508899
$ws.Cells.Item(11, 13).Value = 1.046227376924534
$ws.Cells.Item(11, 14).Value = 1.033593487949488

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.025134234986797
$ws.Cells.Item(12, 4).Value = 1.028168686372347
$ws.Cells.Item(12, 5).Value = 1.034319350475942
$ws.Cells.Item(12, 6).Value = 1.04234862776293
$ws.Cells.Item(12, 9).Value = 1.027130254926545
$ws.Cells.Item(12, 10).Value = 1.031860382753411
$ws.Cells.Item(12, 11).Value = 1.031810434860732
$ws.Cells.Item(12, 12).Value = 1.037937753839226
$ws.Cells.Item(12, 13).Value = 1.045937014506062
$ws.Cells.Item(12, 14).Value = 1.033325742852961

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.025210620346722
$ws.Cells.Item(13, 4).Value = 1.028240319975962
$ws.Cells.Item(13, 5).Value = 1.03438700277544
$ws.Cells.Item(13, 6).Value = 1.042420765732449
$ws.Cells.Item(13, 9).Value = 1.027130625767152
$ws.Cells.Item(13, 10).Value = 1.031917743509937
$ws.Cells.Item(13, 11).Value = 1.031872071156671
$ws.Cells.Item(13, 12).Value = 1.037995493846348
$ws.Cells.Item(13, 13).Value = 1.045999313282914
$ws.Cells.Item(13, 14).Value = 1.033383185068341

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.025460874583963
$ws.Cells.Item(14, 4).Value = 1.028475016621922
$ws.Cells.Item(14, 5).Value = 1.034608647551608
$ws.Cells.Item(14, 6).Value = 1.042657089864312
$ws.Cells.Item(14, 9).Value = 1.027131762869742
$ws.Cells.Item(14, 10).Value = 1.03210565213985
$ws.Cells.Item(14, 11).Value = 1.032073994657188
$ws.Cells.Item(14, 12).Value = 1.038184643316223
$ws.Cells.Item(14, 13).Value = 1.046203382255954
$ws.Cells.Item(14, 14).Value = 1.033571360550061

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.025615043120989
$ws.Cells.Item(15, 4).Value = 1.02861960869128
$ws.Cells.Item(15, 5).Value = 1.034745192692639
$ws.Cells.Item(15, 6).Value = 1.042802665454065
$ws.Cells.Item(15, 9).Value = 1.02713240377623
$ws.Cells.Item(15, 10).Value = 1.032221399735092
$ws.Cells.Item(15, 11).Value = 1.032198381681256
$ws.Cells.Item(15, 12).Value = 1.03830115404886
$ws.Cells.Item(15, 13).Value = 1.046329071880331
$ws.Cells.Item(15, 14).Value = 1.033687272520165

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.02651204957064
$ws.Cells.Item(16, 4).Value = 1.029461013595102
$ws.Cells.Item(16, 5).Value = 1.035539680488125
$ws.Cells.Item(16, 6).Value = 1.043649500840092
$ws.Cells.Item(16, 9).Value = 1.027135222127453
$ws.Cells.Item(16, 10).Value = 1.032894659692443
$ws.Cells.Item(16, 11).Value = 1.032921994293371
$ws.Cells.Item(16, 12).Value = 1.038978834431674
$ws.Cells.Item(16, 13).Value = 1.047059971080831
$ws.Cells.Item(16, 14).Value = 1.03436148858388

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.027074440696917
$ws.Cells.Item(17, 4).Value = 1.029988647510107
$ws.Cells.Item(17, 5).Value = 1.036037813103451
$ws.Cells.Item(17, 6).Value = 1.044180281696327
$ws.Cells.Item(17, 9).Value = 1.027136186270597
$ws.Cells.Item(17, 10).Value = 1.033316593141574
$ws.Cells.Item(17, 11).Value = 1.033375572243873
$ws.Cells.Item(17, 12).Value = 1.03940352164697
$ws.Cells.Item(17, 13).Value = 1.047517859280912
$ws.Cells.Item(17, 14).Value = 1.034784021226905

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.027402376134754
$ws.Cells.Item(18, 4).Value = 1.030296352546838
$ws.Cells.Item(18, 5).Value = 1.036328284577319
$ws.Cells.Item(18, 6).Value = 1.044489728593465
$ws.Cells.Item(18, 9).Value = 1.027136458791933
$ws.Cells.Item(18, 10).Value = 1.033562562605005
$ws.Cells.Item(18, 11).Value = 1.033640020974929
$ws.Cells.Item(18, 12).Value = 1.039651090378077
$ws.Cells.Item(18, 13).Value = 1.047784728046271
$ws.Cells.Item(18, 14).Value = 1.035030339995187

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.027514177258704
$ws.Cells.Item(19, 4).Value = 1.030401262863058
$ws.Cells.Item(19, 5).Value = 1.036427314272066
$ws.Cells.Item(19, 6).Value = 1.044595216928104
$ws.Cells.Item(19, 9).Value = 1.027136502533042
$ws.Cells.Item(19, 10).Value = 1.033646408705495
$ws.Cells.Item(19, 11).Value = 1.033730171693692
$ws.Cells.Item(19, 12).Value = 1.039735480598679
$ws.Cells.Item(19, 13).Value = 1.047875687960345
$ws.Cells.Item(19, 14).Value = 1.035114305166756

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.027014111615058
$ws.Cells.Item(20, 4).Value = 1.029932043102197
$ws.Cells.Item(20, 5).Value = 1.035984376585629
$ws.Cells.Item(20, 6).Value = 1.044123349338629
$ws.Cells.Item(20, 9).Value = 1.027136112803094
$ws.Cells.Item(20, 10).Value = 1.033271337897983
$ws.Cells.Item(20, 11).Value = 1.033326919616019
$ws.Cells.Item(20, 12).Value = 1.039357971667912
$ws.Cells.Item(20, 13).Value = 1.047468753928698
$ws.Cells.Item(20, 14).Value = 1.034738701715677

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.02538718540997
$ws.Cells.Item(21, 4).Value = 1.028405906877354
$ws.Cells.Item(21, 5).Value = 1.034543382347951
$ws.Cells.Item(21, 6).Value = 1.042587504870152
$ws.Cells.Item(21, 9).Value = 1.02713144045722
$ws.Cells.Item(21, 10).Value = 1.032050323805133
$ws.Cells.Item(21, 11).Value = 1.032014538375291
$ws.Cells.Item(21, 12).Value = 1.038128949873089
$ws.Cells.Item(21, 13).Value = 1.046143298204571
$ws.Cells.Item(21, 14).Value = 1.033515953642764

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.024363345214616
$ws.Cells.Item(22, 4).Value = 1.02744583265408
$ws.Cells.Item(22, 5).Value = 1.033636610928284
$ws.Cells.Item(22, 6).Value = 1.041620486689896
$ws.Cells.Item(22, 9).Value = 1.027125895640543
$ws.Cells.Item(22, 10).Value = 1.031281356336901
$ws.Cells.Item(22, 11).Value = 1.031188317273096
$ws.Cells.Item(22, 12).Value = 1.037354886937458
$ws.Cells.Item(22, 13).Value = 1.045308013626121
$ws.Cells.Item(22, 14).Value = 1.032745894152524

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.024906192395689
$ws.Cells.Item(23, 4).Value = 1.027954838383862
$ws.Cells.Item(23, 5).Value = 1.034117381355356
$ws.Cells.Item(23, 6).Value = 1.042133252966053
$ws.Cells.Item(23, 9).Value = 1.027129081982447
$ws.Cells.Item(23, 10).Value = 1.031689122217207
$ws.Cells.Item(23, 11).Value = 1.031626416225203
$ws.Cells.Item(23, 12).Value = 1.037765359680263
$ws.Cells.Item(23, 13).Value = 1.045750996997079
$ws.Cells.Item(23, 14).Value = 1.033154239107152

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.027041372032825
$ws.Cells.Item(24, 4).Value = 1.029957620367633
$ws.Cells.Item(24, 5).Value = 1.036008522498106
$ws.Cells.Item(24, 6).Value = 1.044149075080059
$ws.Cells.Item(24, 9).Value = 1.027136146895652
$ws.Cells.Item(24, 10).Value = 1.033291787218759
$ws.Cells.Item(24, 11).Value = 1.033348904001888
$ws.Cells.Item(24, 12).Value = 1.039378554188485
$ws.Cells.Item(24, 13).Value = 1.047490943171949
$ws.Cells.Item(24, 14).Value = 1.034759180076834

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.029515316260203
$ws.Cells.Item(25, 4).Value = 1.032279602518276
$ws.Cells.Item(25, 5).Value = 1.038199931244027
$ws.Cells.Item(25, 6).Value = 1.046482511823909
$ws.Cells.Item(25, 9).Value = 1.027132972879727
$ws.Cells.Item(25, 10).Value = 1.035146225364231
$ws.Cells.Item(25, 11).Value = 1.035343236666787
$ws.Cells.Item(25, 12).Value = 1.041244937277922
$ws.Cells.Item(25, 13).Value = 1.049501843165034
$ws.Cells.Item(25, 14).Value = 1.03661625173718
